# Commit: "Added Noah Wexler as discussant to attendance sheet"
#
# 2024 - Fall: Nov. 6th seminar (Monique Davis) row - discussant was "TBD",
# now assigned to Noah Wexler. Also a number of "Notes" (column J) cells
# get "JMC" noted (job-market-candidate visits), and the existing
# Thanksgiving note gets ", JMC" appended.

$wb = $excel.ActiveWorkbook

$fall2024 = $wb.Worksheets.Item("2024 - Fall")
$fall2024.Range("F11").Value = "Noah Wexler"

$fall2024.Range("J2").Value = "JMC"
$fall2024.Range("J6").Value = "JMC"
$fall2024.Range("J10").Value = "JMC"
$fall2024.Range("J11").Value = "JMC"
$fall2024.Range("J12").Value = "JMC"
$fall2024.Range("J16").Value = "JMC"
$fall2024.Range("J14").Value = "Thanksgiving Nov. 28th, JMC"

$spring2024 = $wb.Worksheets.Item("2024 - Spring")
$spring2024.Range("J2").Value = "JMC"
$spring2024.Range("J3").Value = "JMC"
$spring2024.Range("J6").Value = "JMC"
$spring2024.Range("J10").Value = "JMC"
$spring2024.Range("J14").Value = "JMC"

$fall2023 = $wb.Worksheets.Item("2023 - Fall")
$fall2023.Range("J2").Value = "JMC"
$fall2023.Range("J4").Value = "JMC"
$fall2023.Range("J11").Value = "JMC"
$fall2023.Range("J14").Value = "JMC"
$fall2023.Range("J15").Value = "JMC"

# Restore the active sheet / selection state to roughly match the author's
# final view (cosmetic, but captured in the diff).
$fall2024.Activate()
$fall2024.Range("E10").Select()
